# Remove duplicate CEDS sectors in Canada scaling mapping
#
# The ceds_sector column (C) on the "map" sheet duplicated information
# already implied by the scaling_sector column (B) on a number of rows.
# This clears those redundant ceds_sector values:
#  - most of them are re-formatted with the same grey "needs lookup" fill
#    already used elsewhere in the sheet (format copied from C34);
#  - row 73 keeps its existing bold style after the value is cleared;
#  - rows 50, 69, 70 and 71 simply lose the stray value (no special
#    formatting was applied there to begin with).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("map")

# Rows whose ceds_sector value gets cleared and re-styled with the grey
# "blank lookup" fill (matching the style already used on C34/C39-C45/C47).
$grayRows = @(23, 25, 26, 27, 52, 53, 54, 58, 59, 62, 64, 65)

$ws.Range("C34").Copy()
foreach ($r in $grayRows) {
    $ws.Range("C$r").PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = 0

foreach ($r in $grayRows) {
    $ws.Range("C$r").ClearContents()
}

# Row 73: clear the value only, keep its existing (bold) style.
$ws.Range("C73").ClearContents()

# Rows 50, 69, 70, 71: clear the stray value entirely (no special style).
$ws.Range("C50").ClearContents()
$ws.Range("C69").ClearContents()
$ws.Range("C70").ClearContents()
$ws.Range("C71").ClearContents()

# Update the sheet's last-saved selection/scroll position to match where
# the editor ended up after making the change.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 55
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C74").Select() | Out-Null
